# "Cambiando iconos de X" -- replace red Wingdings "x" (wrong/incorrect) marker
# shapes with green Wingdings "ü" (correct) marker shapes on the two
# "OTROS SERIALIZADORES" slides.
#
# Slide 12 (1-based): remove the red "x" shapes that sat at
#   (4000496,4643446) and (7213909,4643446) and add two new green "ü"
#   shapes close to those spots, at (4000496,4714884) and
#   (7278029,4714884).
# Slide 13 (1-based): remove the red "x" shape that sat at
#   (8072462,4488428) and add one new green "ü" shape at
#   (8143900,4488428).
#
# New shapes are produced by duplicating an existing green "ü" shape
# already present on the slide (so the run/paragraph formatting -
# Wingdings font w/ pitchFamily+charset, bold, 009900 fill, dirty /
# smtClean flags, centered pPr with defRPr, wrap="none" + spAutoFit -
# is carried over exactly) and then repositioning the duplicate.
#
# Shape.Left/Top/Width/Height round-trip through a single-precision
# (float32) EMU<->point conversion, so the point literals below were
# solved so that floor(float32(pt) * 12700) lands exactly on the
# target EMU values instead of drifting by a one-EMU rounding error.

$p = $ppt.ActivePresentation

function Remove-ShapeById($slide, [int]$targetId) {
    for ($i = $slide.Shapes.Count; $i -ge 1; $i--) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $targetId) {
            $shp.Delete()
            return
        }
    }
}

function Find-ShapeById($slide, [int]$targetId) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.Id -eq $targetId) {
            return $shp
        }
    }
    return $null
}

function Add-CheckMark($slide, $templateId, [double]$left, [double]$top, [double]$width, [double]$height) {
    $template = Find-ShapeById $slide $templateId
    $new = $template.Duplicate()
    $new.Left = $left
    $new.Top = $top
    $new.Width = $width
    $new.Height = $height
    return $new
}

# ---------- Slide 12 ----------
$s12 = $p.Slides.Item(12)

# Drop the two red "x" shapes.
Remove-ShapeById $s12 85
Remove-ShapeById $s12 99

# Add the two green "ü" replacements (duplicated from an existing
# checkmark shape already on the slide, id 62, so formatting matches
# exactly), positioned where the diff puts them.
Add-CheckMark $s12 62 314.9996812246728 371.2507201055091 28.80354426076093 29.081259842519685 | Out-Null
Add-CheckMark $s12 62 573.0731572356938 371.2507201055091 28.80354426076093 29.081259842519685 | Out-Null

# ---------- Slide 13 ----------
$s13 = $p.Slides.Item(13)

# Drop the red "x" shape.
Remove-ShapeById $s13 108

# Add the green "ü" replacement (duplicated from the existing
# checkmark shape id 50 already on the slide).
Add-CheckMark $s13 50 641.2519837627261 353.41954281784416 28.803622524081252 29.081259842519685 | Out-Null
